# Implements the "ShhConnetions" SSH log-query columns on the "Semilla 8"
# sheet (sheet3): adds host/usuario/contraseña ssh columns (I:K) together
# with the server IP and the log-query string, restyles the new/extended
# header row, and makes "Semilla 8" the active sheet/tab (it was previously
# "Semilla 9").

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)   # "Semilla 9"
$ws3 = $wb.Worksheets.Item(3)   # "Semilla 8"

# --- Populate new cells on "Semilla 8" --------------------------------
# Insert values in the same order the target shared-string table uses:
# 10.69.60.119, host ssh, usuario ssh, contraseña ssh, consulta_log
$ws3.Range("I2").Value = "10.69.60.119"
$ws3.Range("I1").Value = "host ssh"
$ws3.Range("J1").Value = "usuario ssh"
$ws3.Range("K1").Value = "contraseña ssh"
$ws3.Range("J2").Value = "consulta_log"
$ws3.Range("K2").Value = "consulta_log"

# --- Restyle header row F1:K1 to match the bold/centered header style ---
# used by A1:E1 (copy format only, keeps existing values/strings intact).
$ws3.Range("A1").Copy()
$ws3.Range("F1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update selections ---------------------------------------------------
$ws2.Activate()
$ws2.Range("B2").Select()

$ws3.Activate()
$ws3.Range("B2").Select()
